$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.384502291679382
$ws.Range("B1").Value = 3.698095560073853
$ws.Range("C1").Value = 3.981334924697876
$ws.Range("D1").Value = 1.696511507034302
$ws.Range("E1").Value = 1.063867688179016
